# Regenerate the "K" (strikeouts) column (G) values for rows 2-40 on the
# active worksheet, replacing the previous "Strike#" derived values with
# the newly computed K values (regen save_data to use K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, indexed by row number (row 2 = first data row).
$newK = @(7, 9, 7, 4, 11, 9, 5, 11, 9, 10, 6, 3, 9, 10, 13, 10, 11, 9, 10, 9, 6, 8, 7, 9, 11, 7, 12, 13, 8, 13, 9, 9, 9, 9, 4, 8, 5, 3, 1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
